$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 602 ("「私たちが愛する言葉」" post entry) entirely.
# This shifts all subsequent rows (603-638) up by one, so the former
# row 603 becomes row 602, ..., and the former row 638 becomes row 637.
$ws.Rows("602:602").Delete()
